# Add data for 2022-04-17 (diff shows update from "through 04-08" to "through 04-09")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update header label to reflect new "through" date
$ws.Name = "Through 2022-04-09"
$ws.Range("A5").Value = "April (through 04-09)"

# Update April row (row 5) values for the affected year columns
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 7
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 32

# Update Total row (row 6) values for the affected year columns
$ws.Range("B6").Value = 72
$ws.Range("C6").Value = 135
$ws.Range("F6").Value = 126
$ws.Range("G6").Value = 218
$ws.Range("H6").Value = 445
$ws.Range("I6").Value = 466
